$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 69.59999999999999
$ws.Range("I4").Value = 74
$ws.Range("J4").Value = 52
$ws.Range("K4").Value = 74
$ws.Range("L4").Value = 52
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = -280
$ws.Range("H5").Value = 67
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 34
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 34
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = -264
$ws.Range("H17").Value = 387340.03
$ws.Range("J17").Value = 394687.9
$ws.Range("L17").Value = 1184063.7
$ws.Range("N17").Value = -1184399.7
$ws.Range("H28").Value = 224700.78
$ws.Range("I28").Value = 252412.88
$ws.Range("K28").Value = 252412.88
$ws.Range("M28").Value = -251927.88
$ws.Range("H32").Value = 2398.75
$ws.Range("I32").Value = 750
$ws.Range("J32").Value = 2948.3333
$ws.Range("K32").Value = 750
$ws.Range("L32").Value = 2948.3333
$ws.Range("M32").Value = -424
$ws.Range("N32").Value = -3600.3333
$ws.Range("H33").Value = 453.6
$ws.Range("I33").Value = 235.1579
$ws.Range("K33").Value = 235.1579
$ws.Range("M33").Value = -6.157900000000012
$ws.Range("H40").Value = 97257.22
$ws.Range("I40").Value = 377042
$ws.Range("J40").Value = 3995.625
$ws.Range("K40").Value = 377042
$ws.Range("L40").Value = 3995.625
$ws.Range("M40").Value = -376867
$ws.Range("N40").Value = -4345.625
$ws.Range("H43").Value = 12479.15
$ws.Range("I43").Value = 11633.167
$ws.Range("J43").Value = 12841.714
$ws.Range("K43").Value = 11633.167
$ws.Range("L43").Value = 12841.714
$ws.Range("M43").Value = -11564.167
$ws.Range("N43").Value = -12979.714
$ws.Range("H51").Value = 3345.7036
$ws.Range("I51").Value = 4325.467
$ws.Range("J51").Value = 2121
$ws.Range("K51").Value = 4325.467
$ws.Range("L51").Value = 2121
$ws.Range("M51").Value = -3841.467
$ws.Range("N51").Value = -3089
$ws.Range("H74").Value = 5805.4443
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 4031.125
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 4031.125
$ws.Range("M74").Value = -19064
$ws.Range("N74").Value = -5903.125
$ws.Range("H76").Value = 8571.429
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 9166.666999999999
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 9166.666999999999
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -9796.666999999999
$ws.Range("H77").Value = 5805.4443
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 4031.125
$ws.Range("K77").Value = 100000
$ws.Range("L77").Value = 20155.625
$ws.Range("M77").Value = -95320
$ws.Range("N77").Value = -29515.625
$ws.Range("H79").Value = 8571.429
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 9166.666999999999
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 9166.666999999999
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -11350.667
$ws.Range("H88").Value = 4324.5
$ws.Range("J88").Value = 4500
$ws.Range("L88").Value = 4500
$ws.Range("N88").Value = -5312
$ws.Range("H91").Value = 4324.5
$ws.Range("J91").Value = 4500
$ws.Range("L91").Value = 4500
$ws.Range("N91").Value = -7308
$ws.Range("H100").Value = 3848.5789
$ws.Range("J100").Value = 4889.421
$ws.Range("L100").Value = 4889.421
$ws.Range("N100").Value = -5971.421
$ws.Range("H111").Value = 1314.6428
$ws.Range("J111").Value = 2453
$ws.Range("L111").Value = 7359
$ws.Range("N111").Value = -13493
$ws.Range("H115").Value = 1415.7142
$ws.Range("I115").Value = 818.3333
$ws.Range("K115").Value = 2454.9999
$ws.Range("M115").Value = -887.9998999999998
$ws.Range("H118").Value = 1066.7693
$ws.Range("J118").Value = 1233
$ws.Range("L118").Value = 3699
$ws.Range("N118").Value = -7013
$ws.Range("H135").Value = 1371.8
$ws.Range("I135").Value = 1371.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 12346.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9811.199999999999
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 452.53333
$ws.Range("I2").Value = 198.07692
$ws.Range("K2").Value = 198.07692
$ws.Range("M2").Value = -85.07692
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H32").Value = 16591148
$ws.Range("I32").Value = 16385799
$ws.Range("K32").Value = 16385799
$ws.Range("M32").Value = -16385512
$ws.Range("H45").Value = 3907.8572
$ws.Range("I45").Value = 3565.9546
$ws.Range("K45").Value = 3565.9546
$ws.Range("M45").Value = -3188.9546
$ws.Range("H61").Value = 3925.2144
$ws.Range("I61").Value = 3904.077
$ws.Range("K61").Value = 3904.077
$ws.Range("M61").Value = -3692.077
$ws.Range("H63").Value = 4966.3335
$ws.Range("I63").Value = 2520
$ws.Range("J63").Value = 6713.7144
$ws.Range("K63").Value = 2520
$ws.Range("L63").Value = 6713.7144
$ws.Range("M63").Value = -1834
$ws.Range("N63").Value = -8085.7144
$ws.Range("H66").Value = 4966.3335
$ws.Range("I66").Value = 2520
$ws.Range("J66").Value = 6713.7144
$ws.Range("K66").Value = 12600
$ws.Range("L66").Value = 33568.572
$ws.Range("M66").Value = -9168
$ws.Range("N66").Value = -40432.572
$ws.Range("H88").Value = 2560.75
$ws.Range("I88").Value = 200
$ws.Range("J88").Value = 2898
$ws.Range("K88").Value = 200
$ws.Range("L88").Value = 2898
$ws.Range("M88").Value = 206
$ws.Range("N88").Value = -3710
$ws.Range("H91").Value = 2560.75
$ws.Range("I91").Value = 200
$ws.Range("J91").Value = 2898
$ws.Range("K91").Value = 200
$ws.Range("L91").Value = 2898
$ws.Range("M91").Value = 1204
$ws.Range("N91").Value = -5706
$ws.Range("H97").Value = 1014.35
$ws.Range("I97").Value = 851.7692
$ws.Range("J97").Value = 1316.2858
$ws.Range("K97").Value = 851.7692
$ws.Range("L97").Value = 1316.2858
$ws.Range("M97").Value = -355.7692
$ws.Range("N97").Value = -2308.2858
$ws.Range("H116").Value = 452.53333
$ws.Range("I116").Value = 198.07692
$ws.Range("K116").Value = 198.07692
$ws.Range("M116").Value = 2095.92308
$ws.Range("H122").Value = 3671.7083
$ws.Range("I122").Value = 3460.0908
$ws.Range("K122").Value = 10380.2724
$ws.Range("M122").Value = -7930.2724
$ws.Range("H132").Value = 5634.6
$ws.Range("I132").Value = 4965.643
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 14896.929
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -12366.929
$ws.Range("N132").Value = -50060
$ws.Range("H136").Value = 3925.2144
$ws.Range("I136").Value = 3904.077
$ws.Range("K136").Value = 11712.231
$ws.Range("M136").Value = -9162.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 452.53333
$ws.Range("I3").Value = 198.07692
$ws.Range("K3").Value = 198.07692
$ws.Range("M3").Value = -84.07692
$ws.Range("H5").Value = 1619.8
$ws.Range("I5").Value = 199.66667
$ws.Range("K5").Value = 199.66667
$ws.Range("M5").Value = -86.66667000000001
$ws.Range("H20").Value = 14894.23
$ws.Range("I20").Value = 20489.223
$ws.Range("J20").Value = 2305.5
$ws.Range("K20").Value = 20489.223
$ws.Range("L20").Value = 2305.5
$ws.Range("M20").Value = -20242.223
$ws.Range("N20").Value = -2799.5
$ws.Range("H99").Value = 2648.2083
$ws.Range("I99").Value = 2455.0952
$ws.Range("K99").Value = 2455.0952
$ws.Range("M99").Value = -957.0952000000002
$ws.Range("H102").Value = 80815
$ws.Range("I102").Value = 20456
$ws.Range("K102").Value = 20456
$ws.Range("M102").Value = -17211
$ws.Range("H107").Value = 1542
$ws.Range("I107").Value = 1466.7878
$ws.Range("J107").Value = 1852.25
$ws.Range("K107").Value = 1466.7878
$ws.Range("L107").Value = 1852.25
$ws.Range("M107").Value = 453.2121999999999
$ws.Range("N107").Value = -5692.25
$ws.Range("H140").Value = 249985.67
$ws.Range("J140").Value = 249985.67
$ws.Range("L140").Value = 249985.67
$ws.Range("N140").Value = -260345.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 3000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = -2713
$ws.Range("H18").Value = 92491.75
$ws.Range("J18").Value = 92491.75
$ws.Range("L18").Value = 92491.75
$ws.Range("N18").Value = -92951.75
$ws.Range("H22").Value = 3595
$ws.Range("I22").Value = 1117.75
$ws.Range("K22").Value = 1117.75
$ws.Range("M22").Value = -767.75
$ws.Range("H28").Value = 57523.5
$ws.Range("J28").Value = 57523.5
$ws.Range("L28").Value = 57523.5
$ws.Range("N28").Value = -58013.5
$ws.Range("H31").Value = 4070.2654
$ws.Range("I31").Value = 2590.6667
$ws.Range("J31").Value = 4276.7207
$ws.Range("K31").Value = 2590.6667
$ws.Range("L31").Value = 4276.7207
$ws.Range("M31").Value = -2295.6667
$ws.Range("N31").Value = -4866.7207
$ws.Range("H34").Value = 4070.2654
$ws.Range("I34").Value = 2590.6667
$ws.Range("J34").Value = 4276.7207
$ws.Range("K34").Value = 2590.6667
$ws.Range("L34").Value = 4276.7207
$ws.Range("M34").Value = -2388.6667
$ws.Range("N34").Value = -4680.7207
$ws.Range("H56").Value = 2500
$ws.Range("J56").Value = 2000
$ws.Range("L56").Value = 2000
$ws.Range("N56").Value = -3690
$ws.Range("H93").Value = 27141.824
$ws.Range("I93").Value = 8036.9165
$ws.Range("J93").Value = 72993.60000000001
$ws.Range("K93").Value = 8036.9165
$ws.Range("L93").Value = 72993.60000000001
$ws.Range("M93").Value = -6164.9165
$ws.Range("N93").Value = -76737.60000000001
$ws.Range("H94").Value = 1112.0625
$ws.Range("I94").Value = 499
$ws.Range("J94").Value = 1253.5385
$ws.Range("K94").Value = 499
$ws.Range("L94").Value = 1253.5385
$ws.Range("M94").Value = -48
$ws.Range("N94").Value = -2155.5385
$ws.Range("H95").Value = 82120.60000000001
$ws.Range("J95").Value = 82120.60000000001
$ws.Range("L95").Value = 82120.60000000001
$ws.Range("N95").Value = -87612.60000000001
$ws.Range("H96").Value = 67324.8
$ws.Range("J96").Value = 67324.8
$ws.Range("L96").Value = 67324.8
$ws.Range("N96").Value = -72816.8
$ws.Range("H102").Value = 188989
$ws.Range("J102").Value = 188989
$ws.Range("L102").Value = 188989
$ws.Range("N102").Value = -193857
$ws.Range("H104").Value = 106990
$ws.Range("J104").Value = 106990
$ws.Range("L104").Value = 106990
$ws.Range("N104").Value = -112232
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H122").Value = 9096405
$ws.Range("I122").Value = 20003996
$ws.Range("K122").Value = 60011988
$ws.Range("M122").Value = -60009538
$ws.Range("H127").Value = 117496.25
$ws.Range("J127").Value = 117496.25
$ws.Range("L127").Value = 117496.25
$ws.Range("N127").Value = -127416.25
$ws.Range("H134").Value = 5205.533
$ws.Range("I134").Value = 5321.273
$ws.Range("J134").Value = 4887.25
$ws.Range("K134").Value = 15963.819
$ws.Range("L134").Value = 14661.75
$ws.Range("M134").Value = -13428.819
$ws.Range("N134").Value = -19731.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 113897270
$ws.Range("I4").Value = 86968136
$ws.Range("K4").Value = 260904408
$ws.Range("M4").Value = -260904296
$ws.Range("H5").Value = 466.4
$ws.Range("I5").Value = 496.875
$ws.Range("J5").Value = 344.5
$ws.Range("K5").Value = 1490.625
$ws.Range("L5").Value = 1033.5
$ws.Range("M5").Value = -1378.625
$ws.Range("N5").Value = -1257.5
$ws.Range("H68").Value = 2137.6052
$ws.Range("J68").Value = 2475.8262
$ws.Range("L68").Value = 7427.4786
$ws.Range("N68").Value = -9049.4786
$ws.Range("H71").Value = 2137.6052
$ws.Range("J71").Value = 2475.8262
$ws.Range("L71").Value = 22282.4358
$ws.Range("N71").Value = -30394.4358
$ws.Range("H98").Value = 2510
$ws.Range("J98").Value = 2588.1667
$ws.Range("L98").Value = 7764.500100000001
$ws.Range("N98").Value = -10760.5001
$ws.Range("H107").Value = 1090.566
$ws.Range("I107").Value = 1036.0769
$ws.Range("J107").Value = 1108.275
$ws.Range("K107").Value = 3108.2307
$ws.Range("L107").Value = 3324.825
$ws.Range("M107").Value = -1188.2307
$ws.Range("N107").Value = -7164.825000000001
$ws.Range("H118").Value = 403
$ws.Range("I118").Value = 403
$ws.Range("K118").Value = 1209
$ws.Range("M118").Value = 34
$ws.Range("H129").Value = 733.1667
$ws.Range("J129").Value = 1499.5
$ws.Range("L129").Value = 4498.5
$ws.Range("N129").Value = -14498.5
$ws.Range("H135").Value = 466.4
$ws.Range("I135").Value = 496.875
$ws.Range("J135").Value = 344.5
$ws.Range("K135").Value = 4471.875
$ws.Range("L135").Value = 3100.5
$ws.Range("M135").Value = -1936.875
$ws.Range("N135").Value = -8170.5
$ws.Range("H139").Value = 2380.6
$ws.Range("I139").Value = 2380.6
$ws.Range("K139").Value = 7141.799999999999
$ws.Range("M139").Value = -2001.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H80").Value = 2423.077
$ws.Range("I80").Value = 2070
$ws.Range("K80").Value = 2070
$ws.Range("M80").Value = -1072
$ws.Range("H83").Value = 2423.077
$ws.Range("I83").Value = 2070
$ws.Range("K83").Value = 10350
$ws.Range("M83").Value = -5358
$ws.Range("H102").Value = 3578.353
$ws.Range("I102").Value = 3388.1538
$ws.Range("J102").Value = 4196.5
$ws.Range("K102").Value = 3388.1538
$ws.Range("L102").Value = 4196.5
$ws.Range("M102").Value = -1766.1538
$ws.Range("N102").Value = -7440.5
$ws.Range("H107").Value = 1016.63635
$ws.Range("I107").Value = 599.2857
$ws.Range("J107").Value = 1747
$ws.Range("K107").Value = 599.2857
$ws.Range("L107").Value = 1747
$ws.Range("M107").Value = 1320.7143
$ws.Range("N107").Value = -5587
$ws.Range("H113").Value = 44903.7
$ws.Range("I113").Value = 8120.5
$ws.Range("J113").Value = 81686.89999999999
$ws.Range("K113").Value = 8120.5
$ws.Range("L113").Value = 81686.89999999999
$ws.Range("M113").Value = -5950.5
$ws.Range("N113").Value = -86026.89999999999
$ws.Range("H122").Value = 4045
$ws.Range("I122").Value = 2803
$ws.Range("K122").Value = 8409
$ws.Range("M122").Value = -5959
$ws.Range("H125").Value = 111323.4
$ws.Range("J125").Value = 111323.4
$ws.Range("L125").Value = 111323.4
$ws.Range("N125").Value = -116243.4
$ws.Range("H126").Value = 4893.4
$ws.Range("I126").Value = 4850.5713
$ws.Range("J126").Value = 4993.3335
$ws.Range("K126").Value = 14551.7139
$ws.Range("L126").Value = 14980.0005
$ws.Range("M126").Value = -12081.7139
$ws.Range("N126").Value = -19920.0005
$ws.Range("H132").Value = 3729.44
$ws.Range("I132").Value = 3660.111
$ws.Range("J132").Value = 3907.7144
$ws.Range("K132").Value = 10980.333
$ws.Range("L132").Value = 11723.1432
$ws.Range("M132").Value = -8450.332999999999
$ws.Range("N132").Value = -16783.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37040036
$ws.Range("I40").Value = 55558324
$ws.Range("J40").Value = 3466
$ws.Range("K40").Value = 55558324
$ws.Range("L40").Value = 3466
$ws.Range("M40").Value = -55558188
$ws.Range("N40").Value = -3738
$ws.Range("H46").Value = 3038.818
$ws.Range("I46").Value = 1723
$ws.Range("K46").Value = 1723
$ws.Range("M46").Value = -1535
$ws.Range("H61").Value = 9062.85
$ws.Range("I61").Value = 9224.923000000001
$ws.Range("J61").Value = 8761.857
$ws.Range("K61").Value = 9224.923000000001
$ws.Range("L61").Value = 8761.857
$ws.Range("M61").Value = -9022.923000000001
$ws.Range("N61").Value = -9165.857
$ws.Range("H93").Value = 66668420
$ws.Range("J93").Value = 2129
$ws.Range("L93").Value = 2129
$ws.Range("N93").Value = -4625
$ws.Range("H113").Value = 9062.85
$ws.Range("I113").Value = 9224.923000000001
$ws.Range("J113").Value = 8761.857
$ws.Range("K113").Value = 9224.923000000001
$ws.Range("L113").Value = 8761.857
$ws.Range("M113").Value = -7054.923000000001
$ws.Range("N113").Value = -13101.857
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5158.3
$ws.Range("I96").Value = 5481.6
$ws.Range("K96").Value = 5481.6
$ws.Range("M96").Value = -4108.6
$ws.Range("H107").Value = 512.44446
$ws.Range("I107").Value = 512.44446
$ws.Range("K107").Value = 1537.33338
$ws.Range("M107").Value = 382.66662
$ws.Range("H113").Value = 635.61536
$ws.Range("I113").Value = 627.5
$ws.Range("J113").Value = 662.6667
$ws.Range("K113").Value = 1882.5
$ws.Range("L113").Value = 1988.0001
$ws.Range("M113").Value = 287.5
$ws.Range("N113").Value = -6328.0001
$ws.Range("H122").Value = 47626404
$ws.Range("I122").Value = 66674030
$ws.Range("K122").Value = 200022090
$ws.Range("M122").Value = -200019640
$ws.Range("H123").Value = 59000
$ws.Range("J123").Value = 59000
$ws.Range("L123").Value = 59000
$ws.Range("N123").Value = -68800
$ws.Range("H124").Value = 75813.164
$ws.Range("J124").Value = 75813.164
$ws.Range("L124").Value = 75813.164
$ws.Range("N124").Value = -85633.164
$ws.Range("H126").Value = 5332
$ws.Range("I126").Value = 4998
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 14994
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -12524
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 3519.9033
$ws.Range("I132").Value = 3454.2856
$ws.Range("J132").Value = 4132.3335
$ws.Range("K132").Value = 10362.8568
$ws.Range("L132").Value = 12397.0005
$ws.Range("M132").Value = -7832.856800000001
$ws.Range("N132").Value = -17457.0005
$ws.Range("H141").Value = 97999.60000000001
$ws.Range("J141").Value = 97999.60000000001
$ws.Range("L141").Value = 97999.60000000001
$ws.Range("N141").Value = -108359.6
